# Auto-generated edit script: update crypto price/volume table per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text/value updates (names, links, percentages, multi-dot prices) ---
$ws.Range("D2").Value = '62.488.50'
$ws.Range("E2").Value = '  +0.97%  '
$ws.Range("D3").Value = '2.433.06'
$ws.Range("E3").Value = '  +1.03%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("E5").Value = '  +0.89%  '
$ws.Range("E6").Value = '  +2.14%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  +0.31%  '
$ws.Range("E9").Value = '  +1.78%  '
$ws.Range("E10").Value = '  +0.52%  '
$ws.Range("E11").Value = '  +1.17%  '
$ws.Range("E12").Value = '  +1.73%  '
$ws.Range("E13").Value = '  +5.17%  '
$ws.Range("E14").Value = '  +4.36%  '
$ws.Range("D15").Value = '2.865.13'
$ws.Range("E15").Value = '  +0.82%  '
$ws.Range("D16").Value = '62.406.98'
$ws.Range("E16").Value = '  +1.00%  '
$ws.Range("D17").Value = '2.447.35'
$ws.Range("E17").Value = '  +1.55%  '
$ws.Range("E18").Value = '  +0.29%  '
$ws.Range("E19").Value = '  +2.13%  '
$ws.Range("E20").Value = '  +0.86%  '
$ws.Range("E21").Value = '  +1.18%  '
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("E23").Value = '  +2.92%  '
$ws.Range("E24").Value = '  +5.09%  '
$ws.Range("E25").Value = '  +3.94%  '
$ws.Range("E26").Value = '  -1.32%  '
$ws.Range("D27").Value = '0.0₃0991'
$ws.Range("E27").Value = '  +6.47%  '
$ws.Range("D28").Value = '2.550.65'
$ws.Range("E28").Value = '  +1.41%  '
$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("E30").Value = '  +3.29%  '
$ws.Range("E31").Value = '  +3.55%  '
$ws.Range("E32").Value = '  -0.74%  '
$ws.Range("E33").Value = '  +0.88%  '
$ws.Range("E34").Value = '  -0.41%  '
$ws.Range("B35").Value = 'FirstDigitalUSD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("E35").Value = '  -0.03%  '
$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("E36").Value = '  +1.65%  '
$ws.Range("E37").Value = '  +0.60%  '
$ws.Range("E38").Value = '  +1.36%  '
$ws.Range("E39").Value = '  -1.96%  '
$ws.Range("E40").Value = '  -2.92%  '
$ws.Range("E41").Value = '  +1.58%  '
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("E43").Value = '  +8.61%  '
$ws.Range("E44").Value = '  +0.26%  '
$ws.Range("E45").Value = '  +2.03%  '
$ws.Range("E46").Value = '  +1.37%  '
$ws.Range("E47").Value = '  +3.08%  '
$ws.Range("E48").Value = '  +1.80%  '
$ws.Range("E49").Value = '  +2.52%  '
$ws.Range("E50").Value = '  +0.46%  '
$ws.Range("E51").Value = '  +4.77%  '

# --- Numeric-looking price strings: force text storage to preserve exact formatting ---
# (set Text number format, assign literal value, then restore the default "Normal" style
#  so no stray cell style is left behind, matching the original un-styled cells)
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '566.20'
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '145.18'
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.110'
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '5.30'
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.354'
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '26.78'
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.0000179'
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '11.22'
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '6.96'
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '323.32'
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '67.22'
$c.Style = "Normal"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '583.69'
$c.Style = "Normal"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '8.42'
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.43'
$c.Style = "Normal"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '1.87'
$c.Style = "Normal"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '4.83'
$c.Style = "Normal"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.381'
$c.Style = "Normal"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '18.72'
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '5.33'
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '147.84'
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '2.42'
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '148.16'
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '3.66'
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.0533'
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '20.44'
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.601'
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.0230'
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.0920'
$c.Style = "Normal"
